$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data rows 2-7 (will be replaced with the new 2-10 range)
$ws.Range("A2:T7").ClearContents()

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Tgfb2"
$ws.Cells.Item(2,3).Value = "Eng"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 1.465615333333333
$ws.Cells.Item(2,8).Value = 4.396846
$ws.Cells.Item(2,9).Value = 0.04672291954663727
$ws.Cells.Item(2,10).Value = 0.04672291954663728
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 114.155417
$ws.Cells.Item(2,14).Value = 342.466251
$ws.Cells.Item(2,15).Value = 0.6835107367845005
$ws.Cells.Item(2,16).Value = 0.6835107367845005
$ws.Cells.Item(2,17).Value = 167.3079295382607
$ws.Cells.Item(2,18).Value = 1505.771365844346
$ws.Cells.Item(2,19).Value = 0.03193561716404498
$ws.Cells.Item(2,20).Value = 0.03193561716404499

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Tgfb2"
$ws.Cells.Item(3,3).Value = "Eng"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 1.465615333333333
$ws.Cells.Item(3,8).Value = 4.396846
$ws.Cells.Item(3,9).Value = 0.04672291954663727
$ws.Cells.Item(3,10).Value = 0.04672291954663728
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 35.924535
$ws.Cells.Item(3,14).Value = 107.773605
$ws.Cells.Item(3,15).Value = 0.2150997826628812
$ws.Cells.Item(3,16).Value = 0.2150997826628812
$ws.Cells.Item(3,17).Value = 52.65154933887
$ws.Cells.Item(3,18).Value = 473.86394404983
$ws.Cells.Item(3,19).Value = 0.01005008983985696
$ws.Cells.Item(3,20).Value = 0.01005008983985696

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Tgfb2"
$ws.Cells.Item(4,3).Value = "Eng"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 1.465615333333333
$ws.Cells.Item(4,8).Value = 4.396846
$ws.Cells.Item(4,9).Value = 0.04672291954663727
$ws.Cells.Item(4,10).Value = 0.04672291954663728
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 16.93339666666667
$ws.Cells.Item(4,14).Value = 50.80019
$ws.Cells.Item(4,15).Value = 0.1013894805526183
$ws.Cells.Item(4,16).Value = 0.1013894805526183
$ws.Cells.Item(4,17).Value = 24.81784580008222
$ws.Cells.Item(4,18).Value = 223.36061220074
$ws.Cells.Item(4,19).Value = 0.004737212542735331
$ws.Cells.Item(4,20).Value = 0.004737212542735332

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Tgfb2"
$ws.Cells.Item(5,3).Value = "Eng"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 18.88237266666667
$ws.Cells.Item(5,8).Value = 56.64711800000001
$ws.Cells.Item(5,9).Value = 0.6019584804341267
$ws.Cells.Item(5,10).Value = 0.6019584804341268
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 114.155417
$ws.Cells.Item(5,14).Value = 342.466251
$ws.Cells.Item(5,15).Value = 0.6835107367845005
$ws.Cells.Item(5,16).Value = 0.6835107367845005
$ws.Cells.Item(5,17).Value = 2155.525125712736
$ws.Cells.Item(5,18).Value = 19399.72613141462
$ws.Cells.Item(5,19).Value = 0.4114450844752083
$ws.Cells.Item(5,20).Value = 0.4114450844752083

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Tgfb2"
$ws.Cells.Item(6,3).Value = "Eng"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 18.88237266666667
$ws.Cells.Item(6,8).Value = 56.64711800000001
$ws.Cells.Item(6,9).Value = 0.6019584804341267
$ws.Cells.Item(6,10).Value = 0.6019584804341268
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 35.924535
$ws.Cells.Item(6,14).Value = 107.773605
$ws.Cells.Item(6,15).Value = 0.2150997826628812
$ws.Cells.Item(6,16).Value = 0.2150997826628812
$ws.Cells.Item(6,17).Value = 678.3404577467101
$ws.Cells.Item(6,18).Value = 6105.06411972039
$ws.Cells.Item(6,19).Value = 0.1294811383134589
$ws.Cells.Item(6,20).Value = 0.1294811383134589

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Tgfb2"
$ws.Cells.Item(7,3).Value = "Eng"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 18.88237266666667
$ws.Cells.Item(7,8).Value = 56.64711800000001
$ws.Cells.Item(7,9).Value = 0.6019584804341267
$ws.Cells.Item(7,10).Value = 0.6019584804341268
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 16.93339666666667
$ws.Cells.Item(7,14).Value = 50.80019
$ws.Cells.Item(7,15).Value = 0.1013894805526183
$ws.Cells.Item(7,16).Value = 0.1013894805526183
$ws.Cells.Item(7,17).Value = 319.7427063724912
$ws.Cells.Item(7,18).Value = 2877.68435735242
$ws.Cells.Item(7,19).Value = 0.06103225764545958
$ws.Cells.Item(7,20).Value = 0.06103225764545959

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Tgfb2"
$ws.Cells.Item(8,3).Value = "Eng"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 11.020243
$ws.Cells.Item(8,8).Value = 33.060729
$ws.Cells.Item(8,9).Value = 0.351318600019236
$ws.Cells.Item(8,10).Value = 0.351318600019236
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 114.155417
$ws.Cells.Item(8,14).Value = 342.466251
$ws.Cells.Item(8,15).Value = 0.6835107367845005
$ws.Cells.Item(8,16).Value = 0.6835107367845005
$ws.Cells.Item(8,17).Value = 1258.020435106331
$ws.Cells.Item(8,18).Value = 11322.18391595698
$ws.Cells.Item(8,19).Value = 0.2401300351452473
$ws.Cells.Item(8,20).Value = 0.2401300351452473

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Tgfb2"
$ws.Cells.Item(9,3).Value = "Eng"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 11.020243
$ws.Cells.Item(9,8).Value = 33.060729
$ws.Cells.Item(9,9).Value = 0.351318600019236
$ws.Cells.Item(9,10).Value = 0.351318600019236
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 35.924535
$ws.Cells.Item(9,14).Value = 107.773605
$ws.Cells.Item(9,15).Value = 0.2150997826628812
$ws.Cells.Item(9,16).Value = 0.2150997826628812
$ws.Cells.Item(9,17).Value = 395.897105362005
$ws.Cells.Item(9,18).Value = 3563.073948258045
$ws.Cells.Item(9,19).Value = 0.07556855450956536
$ws.Cells.Item(9,20).Value = 0.07556855450956537

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Tgfb2"
$ws.Cells.Item(10,3).Value = "Eng"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 11.020243
$ws.Cells.Item(10,8).Value = 33.060729
$ws.Cells.Item(10,9).Value = 0.351318600019236
$ws.Cells.Item(10,10).Value = 0.351318600019236
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 16.93339666666667
$ws.Cells.Item(10,14).Value = 50.80019
$ws.Cells.Item(10,15).Value = 0.1013894805526183
$ws.Cells.Item(10,16).Value = 0.1013894805526183
$ws.Cells.Item(10,17).Value = 186.6101460820567
$ws.Cells.Item(10,18).Value = 1679.49131473851
$ws.Cells.Item(10,19).Value = 0.03562001036442344
$ws.Cells.Item(10,20).Value = 0.03562001036442344

# Update dimension implicitly handled by engine; ensure sheet recalculated